# New crime data collected - update the 78th Precinct weekly CompStat report:
#  - header text (volume/week-covering dates)
#  - weekly/28-day/YTD/2-year crime figures for rows 15-27 (Murder .. Gr. Larceny)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/5/2024  Through  2/11/2024"

# --- Row 15 (Rape) ---------------------------------------------------------
$ws.Range("N15").Value = -83.333333333333

# --- Row 16 (Robbery) -------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 14
$ws.Range("J16").Value = 15
$ws.Range("K16").Value = -6.666666666666
$ws.Range("L16").Value = 55.555555555555
$ws.Range("M16").Value = -30
$ws.Range("N16").Value = -87.272727272727

# --- Row 17 (Fel. Assault) ---------------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 12.5
$ws.Range("I17").Value = 12
$ws.Range("J17").Value = 17
$ws.Range("K17").Value = -29.411764705882
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = -40

# --- Row 18 (Burglary) ---------------------------------------------------------
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -85.714285714285
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -71.428571428571
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = -51.724137931034
$ws.Range("L18").Value = -6.666666666666
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = -86.666666666666

# --- Row 19 (Gr. Larceny) ---------------------------------------------------------
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 14.285714285714
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -40.425531914893
$ws.Range("I19").Value = 46
$ws.Range("J19").Value = 85
$ws.Range("K19").Value = -45.882352941176
$ws.Range("L19").Value = -24.590163934426
$ws.Range("M19").Value = -9.803921568627
$ws.Range("N19").Value = -8

# --- Row 20 (G.L.A.) ---------------------------------------------------------
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 85.714285714285
$ws.Range("I20").Value = 20
$ws.Range("J20").Value = 15
$ws.Range("K20").Value = 33.333333333333
$ws.Range("L20").Value = 122.222222222222
$ws.Range("M20").Value = 233.333333333333
$ws.Range("N20").Value = -86.301369863013

# --- Row 21 (TOTAL) ---------------------------------------------------------
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -19.230769230769
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 93
$ws.Range("H21").Value = -26.881720430107
$ws.Range("I21").Value = 108
$ws.Range("J21").Value = 163
$ws.Range("K21").Value = -33.742331288343
$ws.Range("L21").Value = 1.886792452830
$ws.Range("M21").Value = 14.893617021276
$ws.Range("N21").Value = -75.342465753424

# --- Row 22 (Transit) -------------------------------------------------------
# C22 switches from a numeric 1 to a literal text "0" (same display/style as D22,
# which already stores "0" as text). Writing a bare numeric-looking string gets
# auto-coerced back to a number by the General number format, so format the
# cell as text first, write the value, then restore General formatting/font so
# the cell keeps looking the same as its neighbours.
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C22").NumberFormat = "General"
$ws.Range("C22").Font.Name = "Andale WT"
$ws.Range("C22").Font.Size = 10
$ws.Range("L22").Value = 0

# --- Row 23 (Housing) -------------------------------------------------------
$ws.Range("D23").Value = 2
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -75
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = -25
$ws.Range("L23").Value = 50

# --- Row 24 (Petit Larceny) ---------------------------------------------------------
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -48.484848484848
$ws.Range("F24").Value = 72
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = -45.454545454545
$ws.Range("I24").Value = 113
$ws.Range("J24").Value = 199
$ws.Range("K24").Value = -43.21608040201
$ws.Range("L24").Value = 2.727272727272
$ws.Range("M24").Value = 25.555555555555

# --- Row 25 (Misd. Assault) ---------------------------------------------------------
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 22
$ws.Range("I25").Value = 19
$ws.Range("J25").Value = 33
$ws.Range("K25").Value = -42.424242424242
$ws.Range("L25").Value = -24
$ws.Range("M25").Value = -29.629629629629

# --- Row 27 (Other Sex Crimes) ---------------------------------------------------------
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 14
$ws.Range("K27").Value = 180
$ws.Range("L27").Value = 133.333333333333
